$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so values such as
# "1.006" or "20.226.18" are stored as text, not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "20.226.18"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "1.443.49"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "0.9122"
$ws.Range("E5").Value = "  -8.98%  "
$ws.Range("D6").Value = "278.09"
$ws.Range("E6").Value = "  +3.38%  "
$ws.Range("D7").Value = "0.3657"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").Value = "0.3121"
$ws.Range("E8").Value = "  +3.46%  "
$ws.Range("D9").Value = "39.14"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "1.022"
$ws.Range("E10").Value = "  +6.62%  "
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "5.392"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "17.65"
$ws.Range("E14").Value = "  +7.88%  "
$ws.Range("D15").Value = "6.062"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("E16").Value = "  +2.74%  "
$ws.Range("D17").Value = "1.444.45"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "0.9416"
$ws.Range("E18").Value = "  -6.00%  "
$ws.Range("D19").Value = "0.05637"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "68.12"
$ws.Range("E20").Value = "  -4.14%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "14.43"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.398"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("D24").Value = "2.248"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "20.229.22"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").Value = "2.173"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").Value = "137.84"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("D29").Value = "1.596.13"
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("D30").Value = "110.24"
$ws.Range("E30").Value = "  +3.45%  "
$ws.Range("D31").Value = "3.818"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("D32").Value = "0.8032"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").Value = "4.840"
$ws.Range("E33").Value = "  -6.42%  "
$ws.Range("D34").Value = "0.07694"
$ws.Range("E34").Value = "  +2.17%  "
$ws.Range("D35").Value = "0.05956"
$ws.Range("E35").Value = "  +6.27%  "
$ws.Range("D36").Value = "1.462"
$ws.Range("E36").Value = "  +11.32%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "1.146"
$ws.Range("E37").Value = "  +12.15%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "4.687"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "0.01995"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").Value = "10.16"
$ws.Range("E40").Value = "  +2.85%  "
$ws.Range("D41").Value = "0.9298"
$ws.Range("E41").Value = "  -7.13%  "
$ws.Range("D42").Value = "0.1841"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "7.096"
$ws.Range("E43").Value = "  -14.06%  "
$ws.Range("D44").Value = "3.526"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").Value = "0.5236"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "12.05"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("D47").Value = "118.95"
$ws.Range("E47").Value = "  +10.15%  "
$ws.Range("D48").Value = "0.5142"
$ws.Range("E48").Value = "  +3.39%  "
$ws.Range("D49").Value = "1.759"
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").Value = "0.06338"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("D51").Value = "0.9917"
$ws.Range("E51").Value = "  -1.06%  "
